# Insert a new "USE Violas" sheet, as a duplicate of "USE Cello",
# placed immediately before the "USE Cello" tab, and make it the
# active/selected sheet (matching the commit "Added: Auddict USE Violas").

$wb = $excel.ActiveWorkbook

$cello = $wb.Worksheets.Item("USE Cello")

# Copy "USE Cello" to a new sheet placed right before "USE Cello" itself.
$cello.Copy($cello)

# The newly inserted copy takes the default name "USE Cello (2)".
$violas = $wb.Worksheets.Item("USE Cello (2)")
$violas.Name = "USE Violas"

# Make the new sheet the active tab/selection, matching the target workbook
# (activeTab points at "USE Violas" and its sheetView has tabSelected="1").
$violas.Activate()
$violas.Range("A2").Select()
